$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: ALC
$ws.Range("H17").Value = 480
$ws.Range("J17").Value = 480
$ws.Range("L17").Value = 1440
$ws.Range("N17").Value = -1776

# Row 40: ALC
$ws.Range("H40").Value = 2649.8333
$ws.Range("I40").Value = 2499
$ws.Range("J40").Value = 2680
$ws.Range("K40").Value = 2499
$ws.Range("L40").Value = 2680
$ws.Range("M40").Value = -2324
$ws.Range("N40").Value = -3030

# Row 51: ALC
$ws.Range("H51").Value = 5310
$ws.Range("I51").Value = 6634
$ws.Range("J51").Value = 5111.4
$ws.Range("K51").Value = 6634
$ws.Range("L51").Value = 5111.4
$ws.Range("M51").Value = -6150
$ws.Range("N51").Value = -6079.4

# Row 64: ALC
$ws.Range("H64").Value = 5428.5713
$ws.Range("J64").Value = 6000
$ws.Range("L64").Value = 6000
$ws.Range("N64").Value = -6496

# Row 67: ALC
$ws.Range("H67").Value = 5428.5713
$ws.Range("J67").Value = 6000
$ws.Range("L67").Value = 6000
$ws.Range("N67").Value = -7716

# Row 70: ALC
$ws.Range("H70").Value = 682206.25
$ws.Range("J70").Value = 2369.4
$ws.Range("L70").Value = 7108.200000000001
$ws.Range("N70").Value = -7648.200000000001

# Row 73: ALC
$ws.Range("H73").Value = 682206.25
$ws.Range("J73").Value = 2369.4
$ws.Range("L73").Value = 7108.200000000001
$ws.Range("N73").Value = -8980.200000000001

# Row 100: ALC
$ws.Range("H100").Value = 1589
$ws.Range("I100").Value = 1436.5
$ws.Range("K100").Value = 1436.5
$ws.Range("M100").Value = -895.5

# Row 112: ALC
$ws.Range("H112").Value = 60358.47
$ws.Range("I112").Value = 92243.91
$ws.Range("K112").Value = 276731.73
$ws.Range("M112").Value = -275623.73

# Row 116: ALC
$ws.Range("H116").Value = 83492760
$ws.Range("I116").Value = 50222000
$ws.Range("J116").Value = 166669660
$ws.Range("K116").Value = 50222000
$ws.Range("L116").Value = 166669660
$ws.Range("M116").Value = -50218558
$ws.Range("N116").Value = -166676544

# Row 135: ALC
$ws.Range("H135").Value = 38464948
$ws.Range("I135").Value = 43481850
$ws.Range("J135").Value = 2065
$ws.Range("K135").Value = 391336650
$ws.Range("L135").Value = 18585
$ws.Range("M135").Value = -391334115
$ws.Range("N135").Value = -23655

$ws = $wb.Worksheets.Item("ARM")
# Row 4: ARM
$ws.Range("H4").Value = 7243.357
$ws.Range("I4").Value = 148.75
$ws.Range("K4").Value = 148.75
$ws.Range("M4").Value = -32.75

# Row 45: ARM
$ws.Range("H45").Value = 1881.8823
$ws.Range("I45").Value = 1733.4667
$ws.Range("K45").Value = 1733.4667
$ws.Range("M45").Value = -1356.4667

# Row 110: ARM
$ws.Range("H110").Value = 91003410
$ws.Range("I110").Value = 125065940
$ws.Range("J110").Value = 169983
$ws.Range("K110").Value = 125065940
$ws.Range("L110").Value = 169983
$ws.Range("M110").Value = -125063895
$ws.Range("N110").Value = -174073

# Row 122: ARM
$ws.Range("H122").Value = 14495952
$ws.Range("I122").Value = 18521496
$ws.Range("K122").Value = 55564488
$ws.Range("M122").Value = -55562038

$ws = $wb.Worksheets.Item("BSM")
# Row 105: BSM
$ws.Range("H105").Value = 2032.8334
$ws.Range("I105").Value = 1910.4445
$ws.Range("K105").Value = 1910.4445
$ws.Range("M105").Value = -163.4445000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 16: CRP
$ws.Range("H16").Value = 1848.2222
$ws.Range("I16").Value = 1805.5714
$ws.Range("J16").Value = 1997.5
$ws.Range("K16").Value = 1805.5714
$ws.Range("L16").Value = 1997.5
$ws.Range("M16").Value = -1518.5714
$ws.Range("N16").Value = -2571.5

# Row 31: CRP
$ws.Range("H31").Value = 3411.6897
$ws.Range("I31").Value = 1958.7778
$ws.Range("K31").Value = 1958.7778
$ws.Range("M31").Value = -1663.7778

# Row 34: CRP
$ws.Range("H34").Value = 3411.6897
$ws.Range("I34").Value = 1958.7778
$ws.Range("K34").Value = 1958.7778
$ws.Range("M34").Value = -1756.7778

# Row 58: CRP
$ws.Range("H58").Value = 2019.9445
$ws.Range("I58").Value = 2014.5454
$ws.Range("J58").Value = 2028.4286
$ws.Range("K58").Value = 2014.5454
$ws.Range("L58").Value = 2028.4286
$ws.Range("M58").Value = -1811.5454
$ws.Range("N58").Value = -2434.4286

# Row 105: CRP
$ws.Range("H105").Value = 1770
$ws.Range("I105").Value = 1248.75
$ws.Range("K105").Value = 1248.75
$ws.Range("M105").Value = 498.25

# Row 107: CRP
$ws.Range("H107").Value = 4066.3572
$ws.Range("I107").Value = 2929.9092
$ws.Range("K107").Value = 2929.9092
$ws.Range("M107").Value = -1009.9092

# Row 113: CRP
$ws.Range("H113").Value = 1848.2222
$ws.Range("I113").Value = 1805.5714
$ws.Range("J113").Value = 1997.5
$ws.Range("K113").Value = 1805.5714
$ws.Range("L113").Value = 1997.5
$ws.Range("M113").Value = 364.4286
$ws.Range("N113").Value = -6337.5

# Row 122: CRP
$ws.Range("H122").Value = 3348.4
$ws.Range("I122").Value = 2712.7144
$ws.Range("K122").Value = 8138.1432
$ws.Range("M122").Value = -5688.1432

# Row 134: CRP
$ws.Range("H134").Value = 2389.0588
$ws.Range("I134").Value = 946.25
$ws.Range("K134").Value = 2838.75
$ws.Range("M134").Value = -303.75

# Row 136: CRP
$ws.Range("H136").Value = 2019.9445
$ws.Range("I136").Value = 2014.5454
$ws.Range("J136").Value = 2028.4286
$ws.Range("K136").Value = 6043.6362
$ws.Range("L136").Value = 6085.2858
$ws.Range("M136").Value = -3493.6362
$ws.Range("N136").Value = -11185.2858

$ws = $wb.Worksheets.Item("CUL")
# Row 5: CUL
$ws.Range("H5").Value = 1074.3334
$ws.Range("I5").Value = 1074.3334
$ws.Range("K5").Value = 3223.0002
$ws.Range("M5").Value = -3111.0002

# Row 9: CUL
$ws.Range("H9").Value = 5116
$ws.Range("I9").Value = 7223.1113
$ws.Range("J9").Value = 375
$ws.Range("K9").Value = 21669.3339
$ws.Range("L9").Value = 1125
$ws.Range("M9").Value = -21445.3339
$ws.Range("N9").Value = -1573

# Row 55: CUL
$ws.Range("H55").Value = 483.36
$ws.Range("I55").Value = 469.7
$ws.Range("J55").Value = 538
$ws.Range("K55").Value = 1409.1
$ws.Range("L55").Value = 1614
$ws.Range("M55").Value = -1232.1
$ws.Range("N55").Value = -1968

# Row 135: CUL
$ws.Range("H135").Value = 1074.3334
$ws.Range("I135").Value = 1074.3334
$ws.Range("K135").Value = 9669.000599999999
$ws.Range("M135").Value = -7134.000599999999

$ws = $wb.Worksheets.Item("GSM")
# Row 15: GSM
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null

# Row 81: GSM
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null

# Row 84: GSM
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null

# Row 96: GSM
$ws.Range("H96").Value = 50000.668
$ws.Range("J96").Value = 50000.668
$ws.Range("L96").Value = 50000.668
$ws.Range("N96").Value = -55492.668

# Row 126: GSM
$ws.Range("H126").Value = 52006
$ws.Range("I126").Value = 100012
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 300036
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -297566
$ws.Range("N126").Value = -16940

# Row 132: GSM
$ws.Range("H132").Value = 5075
$ws.Range("I132").Value = 4712.1665
$ws.Range("K132").Value = 14136.4995
$ws.Range("M132").Value = -11606.4995

$ws = $wb.Worksheets.Item("LTW")
# Row 16: LTW
$ws.Range("H16").Value = 3565.4119
$ws.Range("I16").Value = 3473
$ws.Range("K16").Value = 3473
$ws.Range("M16").Value = -3303

# Row 46: LTW
$ws.Range("H46").Value = 1508.2858
$ws.Range("I46").Value = 1101.6666
$ws.Range("J46").Value = 2240.2
$ws.Range("K46").Value = 1101.6666
$ws.Range("L46").Value = 2240.2
$ws.Range("M46").Value = -913.6666
$ws.Range("N46").Value = -2616.2

# Row 55: LTW
$ws.Range("H55").Value = 898.2727
$ws.Range("I55").Value = 788
$ws.Range("K55").Value = 788
$ws.Range("M55").Value = -615

# Row 61: LTW
$ws.Range("H61").Value = 1952
$ws.Range("I61").Value = 1404
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 1404
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -1202
$ws.Range("N61").Value = -2904

# Row 93: LTW
$ws.Range("H93").Value = 1571.5217
$ws.Range("I93").Value = 1459.2858
$ws.Range("K93").Value = 1459.2858
$ws.Range("M93").Value = -211.2858000000001

# Row 113: LTW
$ws.Range("H113").Value = 1952
$ws.Range("I113").Value = 1404
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1404
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 766
$ws.Range("N113").Value = -6840

# Row 132: LTW
$ws.Range("H132").Value = 4604.727
$ws.Range("I132").Value = 2117.1365
$ws.Range("K132").Value = 6351.4095
$ws.Range("M132").Value = -3821.4095

$ws = $wb.Worksheets.Item("WVR")
# Row 122: WVR
$ws.Range("H122").Value = 1796.75
$ws.Range("I122").Value = 1596.0769
$ws.Range("K122").Value = 4788.2307
$ws.Range("M122").Value = -2338.2307

# Row 124: WVR
$ws.Range("H124").Value = 100000
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820

# Row 125: WVR
$ws.Range("H125").Value = 35700
$ws.Range("J125").Value = 35700
$ws.Range("L125").Value = 35700
$ws.Range("N125").Value = -45540

# Row 132: WVR
$ws.Range("H132").Value = 4002.75
$ws.Range("I132").Value = 3818.087
$ws.Range("K132").Value = 11454.261
$ws.Range("M132").Value = -8924.261

# Row 136: WVR
$ws.Range("H136").Value = 4181.15
$ws.Range("I136").Value = 1477
$ws.Range("K136").Value = 4431
$ws.Range("M136").Value = -1881
